# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2304   (the "source"/older format-version columns)
#   *_new -> *_FV2310   (the "target"/newer format-version columns)
# then wrap the data range in a native Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row (row 1) cells ------------------------------
# Row 1 holds the column captions used both as plain cell values and later
# as the Excel Table's column names.
$lastCol = 21   # columns A..U
for ($i = 1; $i -le $lastCol; $i++) {
    $cell = $ws.Cells.Item(1, $i)
    $v = $cell.Value()
    if ($v -like "*_old") {
        $cell.Value = ($v -replace "_old$", "_FV2304")
    } elseif ($v -like "*_new") {
        $cell.Value = ($v -replace "_new$", "_FV2310")
    }
}

# --- 2) Turn the data range into a real Excel Table (ListObject) ---------
# This produces xl/tables/table1.xml, registers it in
# xl/worksheets/_rels/sheet1.xml.rels and adds the <tableParts> entry to
# the worksheet, with column names picked up from the renamed headers.
$dataRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ---------------------------------------------
$win = $excel.ActiveWindow
$ws.Range("A2").Select() | Out-Null
$win.FreezePanes = $true
